# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are text (not numbers, since they use "." as
# a thousands separator as well as a decimal point for some coins), so we
# force the NumberFormat to Text before writing, then restore the default
# "Normal" style so no formatting residue is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.912.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.814.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4648"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8687"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.811.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.359"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07098"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.495"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.004"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008692"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.923.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.290"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.029.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.896"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.246"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08905"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7535"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.158"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.475"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.906"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.086"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05277"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01946"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.236"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5299"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.307"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1650"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.410"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4865"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.659"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("E51").Value = "  -0.08%  "
